$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1966.7142
$ws.Range("J17").Value = 1969.5
$ws.Range("L17").Value = 5908.5
$ws.Range("N17").Value = -6244.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 299.75
$ws.Range("I29").Value = 299.75
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 899.25
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -618.25
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 114.5625
$ws.Range("I33").Value = 114.5625
$ws.Range("K33").Value = 114.5625
$ws.Range("M33").Value = 114.4375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 496.6
$ws.Range("J38").Value = 1999
$ws.Range("L38").Value = 5997
$ws.Range("N38").Value = -6741

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 6199.75
$ws.Range("J41").Value = 6199.75
$ws.Range("L41").Value = 6199.75
$ws.Range("N41").Value = -7079.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 79.5
$ws.Range("I53").Value = 76.71429000000001
$ws.Range("J53").Value = 99
$ws.Range("K53").Value = 76.71429000000001
$ws.Range("L53").Value = 99
$ws.Range("M53").Value = 560.28571
$ws.Range("N53").Value = -1373

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1946.375
$ws.Range("I92").Value = 1514.8
$ws.Range("K92").Value = 1514.8
$ws.Range("M92").Value = -266.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 3848.6
$ws.Range("J97").Value = 3848.6
$ws.Range("L97").Value = 11545.8
$ws.Range("N97").Value = -12537.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1968.6666
$ws.Range("I61").Value = 1968.6666
$ws.Range("K61").Value = 1968.6666
$ws.Range("M61").Value = -1756.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 11336.053
$ws.Range("I122").Value = 11099.214
$ws.Range("J122").Value = 11999.2
$ws.Range("K122").Value = 33297.642
$ws.Range("L122").Value = 35997.60000000001
$ws.Range("M122").Value = -30847.642
$ws.Range("N122").Value = -40897.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1870.7333
$ws.Range("I132").Value = 1870.7333
$ws.Range("K132").Value = 5612.199900000001
$ws.Range("M132").Value = -3082.199900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1968.6666
$ws.Range("I136").Value = 1968.6666
$ws.Range("K136").Value = 5905.9998
$ws.Range("M136").Value = -3355.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7333
$ws.Range("I105").Value = 7333
$ws.Range("K105").Value = 7333
$ws.Range("M105").Value = -5586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5309.4546
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 5000
$ws.Range("N31").Value = -5590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5309.4546
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5404

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 9000000
$ws.Range("I10").Value = 9000000
$ws.Range("K10").Value = 9000000
$ws.Range("M10").Value = -8999831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 100000
$ws.Range("I33").Value = 100000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 100000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -99748
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4228
$ws.Range("I80").Value = 2688
$ws.Range("J80").Value = 4998
$ws.Range("K80").Value = 2688
$ws.Range("L80").Value = 4998
$ws.Range("M80").Value = -1690
$ws.Range("N80").Value = -6994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4228
$ws.Range("I83").Value = 2688
$ws.Range("J83").Value = 4998
$ws.Range("K83").Value = 13440
$ws.Range("L83").Value = 24990
$ws.Range("M83").Value = -8448
$ws.Range("N83").Value = -34974

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3055.2307
$ws.Range("I102").Value = 2940
$ws.Range("K102").Value = 2940
$ws.Range("M102").Value = -1318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7404.5557
$ws.Range("I122").Value = 7829.125
$ws.Range("K122").Value = 23487.375
$ws.Range("M122").Value = -21037.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3678
$ws.Range("I132").Value = 3880.3333
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 11640.9999
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -9110.999899999999
$ws.Range("N132").Value = -8810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 450
$ws.Range("I14").Value = 450
$ws.Range("K14").Value = 450
$ws.Range("M14").Value = -278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 30005
$ws.Range("I20").Value = 30005
$ws.Range("K20").Value = 30005
$ws.Range("M20").Value = -29779

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 16503
$ws.Range("I21").Value = 16503
$ws.Range("K21").Value = 16503
$ws.Range("M21").Value = -16329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 1199.6666
$ws.Range("L22").Value = 1199.6666
$ws.Range("N22").Value = -1789.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 30006
$ws.Range("I24").Value = 30006
$ws.Range("K24").Value = 30006
$ws.Range("M24").Value = -29663

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J27").Value = 1199.6666
$ws.Range("L27").Value = 1199.6666
$ws.Range("N27").Value = -1413.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10000
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 7750
$ws.Range("I12").Value = 6500
$ws.Range("J12").Value = 9000
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = -6358
$ws.Range("N12").Value = -9284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 30006
$ws.Range("I15").Value = 30006
$ws.Range("K15").Value = 30006
$ws.Range("M15").Value = -29718

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 30005
$ws.Range("I19").Value = 30005
$ws.Range("K19").Value = 30005
$ws.Range("M19").Value = -29831

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 24262
$ws.Range("I37").Value = 13525.5
$ws.Range("J37").Value = 34998.5
$ws.Range("K37").Value = 13525.5
$ws.Range("L37").Value = 34998.5
$ws.Range("M37").Value = -13322.5
$ws.Range("N37").Value = -35404.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3221.7778
$ws.Range("I122").Value = 2998.6
$ws.Range("K122").Value = 8995.799999999999
$ws.Range("M122").Value = -6545.799999999999
